# Auto-generated script applying 2023-08-23 violent-crime daily data update
# across Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 7280
$ws.Range('J2').Value = 4849
$ws.Range('J3').Value = 5139
$ws.Range('B4').Value = 1682
$ws.Range('I4').Value = 1774
$ws.Range('J4').Value = 1144
$ws.Range('J5').Value = 409
$ws.Range('J6').Value = 6366
$ws.Range('B7').Value = 23314
$ws.Range('I7').Value = 26225
$ws.Range('J7').Value = 17907

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 204
$ws.Range('J3').Value = 284
$ws.Range('J6').Value = 222
$ws.Range('J7').Value = 783

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J6').Value = 107
$ws.Range('J7').Value = 385

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J2').Value = 76
$ws.Range('J3').Value = 102
$ws.Range('J6').Value = 73
$ws.Range('J7').Value = 270

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('J3').Value = 25
$ws.Range('J6').Value = 19
$ws.Range('J7').Value = 69

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 171
$ws.Range('J3').Value = 268
$ws.Range('J4').Value = 52
$ws.Range('J5').Value = 20
$ws.Range('J6').Value = 184
$ws.Range('J7').Value = 695

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J2').Value = 61
$ws.Range('J7').Value = 163

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 140
$ws.Range('J5').Value = 56
$ws.Range('J6').Value = 134
$ws.Range('J7').Value = 517
$ws.Range('J8').Value = 1142
$ws.Range('J10').Value = 120
$ws.Range('J11').Value = 273
$ws.Range('J15').Value = 194
$ws.Range('J16').Value = 68
$ws.Range('J18').Value = 156
$ws.Range('J19').Value = 524
$ws.Range('J20').Value = 375
$ws.Range('J21').Value = 46
$ws.Range('J23').Value = 172
$ws.Range('J27').Value = 97
$ws.Range('J29').Value = 1015
$ws.Range('J30').Value = 69
$ws.Range('J31').Value = 163
$ws.Range('J33').Value = 813
$ws.Range('J36').Value = 250
$ws.Range('J37').Value = 558
$ws.Range('J42').Value = 728
$ws.Range('J47').Value = 138
$ws.Range('J48').Value = 200
$ws.Range('J51').Value = 224
$ws.Range('I63').Value = 236
$ws.Range('J63').Value = 67
$ws.Range('J65').Value = 472
$ws.Range('J67').Value = 695
$ws.Range('J68').Value = 33
$ws.Range('J73').Value = 169
$ws.Range('J77').Value = 138
$ws.Range('J83').Value = 385
$ws.Range('J85').Value = 783
$ws.Range('J89').Value = 226
$ws.Range('B91').Value = 311
$ws.Range('J91').Value = 199
$ws.Range('J93').Value = 76
$ws.Range('J94').Value = 170
$ws.Range('J96').Value = 214
$ws.Range('J97').Value = 140
$ws.Range('J99').Value = 270
$ws.Range('B101').Value = 23314
$ws.Range('I101').Value = 26225
$ws.Range('J101').Value = 17907

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J5').Value = 23
$ws.Range('J7').Value = 558

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J6').Value = 281
$ws.Range('J7').Value = 813

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J6').Value = 167
$ws.Range('J7').Value = 472

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 305
$ws.Range('J3').Value = 351
$ws.Range('J6').Value = 261
$ws.Range('J7').Value = 1015

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 128
$ws.Range('J3').Value = 150
$ws.Range('J6').Value = 196
$ws.Range('J7').Value = 524

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J2').Value = 31
$ws.Range('J7').Value = 200

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 164
$ws.Range('J4').Value = 19
$ws.Range('J7').Value = 517

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J2').Value = 40
$ws.Range('J7').Value = 134

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 159
$ws.Range('J6').Value = 372
$ws.Range('J7').Value = 728

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J6').Value = 67
$ws.Range('J7').Value = 120

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J4').Value = 19
$ws.Range('J6').Value = 42
$ws.Range('J7').Value = 172

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J4').Value = 26
$ws.Range('J7').Value = 226

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J3').Value = 83
$ws.Range('B4').Value = 11
$ws.Range('J6').Value = 39
$ws.Range('B7').Value = 311
$ws.Range('J7').Value = 199

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('J6').Value = 28
$ws.Range('J7').Value = 46

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J3').Value = 126
$ws.Range('J7').Value = 375

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J2').Value = 41
$ws.Range('J7').Value = 156

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J3').Value = 77
$ws.Range('J5').Value = 4
$ws.Range('J6').Value = 71
$ws.Range('J7').Value = 250

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('J6').Value = 28
$ws.Range('J7').Value = 76

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J6').Value = 104
$ws.Range('J7').Value = 273

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J3').Value = 32
$ws.Range('J7').Value = 170

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('J3').Value = 39
$ws.Range('J7').Value = 138

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J6').Value = 82
$ws.Range('J7').Value = 194

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J2').Value = 59
$ws.Range('J6').Value = 51
$ws.Range('J7').Value = 169

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J2').Value = 39
$ws.Range('J6').Value = 53
$ws.Range('J7').Value = 140

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J2').Value = 24
$ws.Range('J7').Value = 140

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('J6').Value = 26
$ws.Range('J7').Value = 56

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J2').Value = 27
$ws.Range('J7').Value = 97

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J3').Value = 63
$ws.Range('J7').Value = 224

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('J4').Value = 3
$ws.Range('J7').Value = 33

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J2').Value = 66
$ws.Range('J7').Value = 214

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J3').Value = 50
$ws.Range('J7').Value = 138

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 321
$ws.Range('J3').Value = 344
$ws.Range('J6').Value = 382
$ws.Range('J7').Value = 1142

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 68
